# Trade #124 closed at 2026-02-17 16:07:12 - unknown UNKNOWN +0.000%
#
# Updates the Summary and Strategy Status roll-up figures to reflect the
# newly closed MarketMaking trade, and appends the new trade row (#124,
# 1-based spreadsheet row 125) to both the "All Trades" and "MarketMaking"
# sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.48   # Current Capital
$summary.Range("B4").Value = -0.53     # Total P&L $
$summary.Range("B5").Value = -0.09     # Total P&L %
$summary.Range("B6").Value = 124       # Total Trades
$summary.Range("B8").Value = 61        # Losing Trades
$summary.Range("B9").Value = 37.9      # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.48      # Capital
$status.Range("D4").Value = 124        # Trades
$status.Range("E4").Value = -0.53      # P&L $
$status.Range("F4").Value = -0.52      # P&L %
$status.Range("G4").Value = 37.9       # Win Rate %

# ---------------------------------------------------------------------
# New trade row (Trade #124 -> sheet row 125) shared by "All Trades" and
# "MarketMaking" sheets.
# ---------------------------------------------------------------------
$newRow = @(124, "2026-02-17", "16:07:05", "MarketMaking", "UP", 0.39, 0.24558, "CLOSED", -37.0309, -0.14, 99.48, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column B holds a plain "YYYY-MM-DD" text label (matches every other
    # row in the sheet, which is stored as literal text, not a date
    # serial). Force text formatting first so Excel doesn't auto-convert
    # the string into a date value on assignment.
    $ws.Cells.Item(125, 2).NumberFormat = "@"

    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item(125, $col).Value = $newRow[$col - 1]
    }
}
